# This workbook's "Artfynd" sheet rows 3-18 are a single logical block that
# got reshuffled: each row keeps all of its populated columns (A:AY - id,
# taxon info, coordinates, observers, etc.) together, it just moves to a
# different row position. We implement that as a permutation: stash the
# current rows 3-18 (columns A:AY, which is where all the data lives) in
# scratch rows far below the used range, then copy them back into their new
# homes in the right order, then wipe the scratch rows.
#
# Note: we deliberately use bounded ranges ("A<row>:AY<row>") rather than
# whole-row objects (Rows.Item) for both Copy and Clear - whole-row
# operations here are drastically slower and also balloon the sheet's used
# range out to column XFD.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 3
$lastRow = 18
$firstCol = "A"
$lastCol = "AY"
$scratchOffset = 200   # scratch block starts at row 203, well clear of used range (A1:AY18)

# Map: target row -> source row (both in original 3..18 numbering)
$targetToSource = @{
    3  = 14
    4  = 7
    5  = 6
    6  = 12
    7  = 10
    8  = 18
    9  = 5
    10 = 15
    11 = 3
    12 = 17
    13 = 9
    14 = 8
    15 = 16
    16 = 4
    17 = 11
    18 = 13
}

# Step 1: copy each original row's data (A:AY) down into its own scratch row
# so that later overwrites of rows 3..18 can't clobber data we still need.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $srcRange = $ws.Range($firstCol + $r + ":" + $lastCol + $r)
    $dstRowNum = $r + $scratchOffset
    $dstRange = $ws.Range($firstCol + $dstRowNum + ":" + $lastCol + $dstRowNum)
    $srcRange.Copy($dstRange)
}

# Step 2: copy from the scratch rows into the final target rows, in the
# order dictated by the permutation.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $srcRow = $targetToSource[$r]
    $scratchRowNum = $srcRow + $scratchOffset
    $scratchRange = $ws.Range($firstCol + $scratchRowNum + ":" + $lastCol + $scratchRowNum)
    $destRange = $ws.Range($firstCol + $r + ":" + $lastCol + $r)
    $scratchRange.Copy($destRange)
}

# Step 3: clear the scratch rows so the saved workbook's used range matches
# the original extent.
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $scratchRowNum = $r + $scratchOffset
    $ws.Range($firstCol + $scratchRowNum + ":" + $lastCol + $scratchRowNum).Clear()
}

Write-Output "row permutation applied"
